$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings) ---
$ws.Range("A8").Value = "Volume 31   Number  6"
$ws.Range("C9").Value = "Report Covering the Week  2/5/2024  Through  2/11/2024"

# --- Crime-data table updates (Week 2/5/2024 - 2/11/2024) ---
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = -33.333333333333
$ws.Range("I16").Value = 7
$ws.Range("J16").Value = 11
$ws.Range("K16").Value = -36.363636363636
$ws.Range("L16").Value = 16.666666666666
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -25
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = -35
$ws.Range("I17").Value = 17
$ws.Range("J17").Value = 29
$ws.Range("K17").Value = -41.379310344827
$ws.Range("L17").Value = -15
$ws.Range("C18").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("D18").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E18").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("C19").Value = 5
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 24
$ws.Range("G19").Value = 19
$ws.Range("H19").Value = 26.315789473684
$ws.Range("I19").Value = 39
$ws.Range("J19").Value = 34
$ws.Range("K19").Value = 14.705882352941
$ws.Range("L19").Value = -29.090909090909
$ws.Range("D20").Value = 1
$ws.Range("F18").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").Value = -100
$ws.Range("K18").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = -25
$ws.Range("I20").Value = 3
$ws.Range("J20").Value = 7
$ws.Range("K20").Value = -57.142857142857
$ws.Range("L20").Value = -78.571428571428
$ws.Range("C21").Value = 9
$ws.Range("D21").Value = 13
$ws.Range("E21").Value = -30.769230769230
$ws.Range("F21").Value = 52
$ws.Range("G21").Value = 57
$ws.Range("H21").Value = -8.771929824561
$ws.Range("I21").Value = 73
$ws.Range("J21").Value = 89
$ws.Range("K21").Value = -17.977528089887
$ws.Range("L21").Value = -30.476190476190
$ws.Range("D23").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E23").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 0
$ws.Range("C24").Value = 34
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = 41.666666666666
$ws.Range("F24").Value = 122
$ws.Range("G24").Value = 83
$ws.Range("H24").Value = 46.987951807228
$ws.Range("I24").Value = 184
$ws.Range("J24").Value = 142
$ws.Range("K24").Value = 29.577464788732
$ws.Range("L24").Value = 38.345864661654
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 37
$ws.Range("G25").Value = 39
$ws.Range("H25").Value = -5.128205128205
$ws.Range("I25").Value = 56
$ws.Range("J25").Value = 58
$ws.Range("K25").Value = -3.448275862068
$ws.Range("L25").Value = 16.666666666666
$ws.Range("D26").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("C27").Value = 2
$ws.Range("E27").Value = 100
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 6
$ws.Range("J27").Value = 8
$ws.Range("K27").Value = -25
$ws.Range("L27").Value = 20
